# Estonia Meistriliiga - base update (20-06-2024 20:11)
# Two pairs of fixture rows were recorded swapped in the source feed and
# are corrected here by swapping the match data between each row pair
# (match id, teams, score, odds, etc.) while leaving the row's sequence
# number (column A), league name (column C) and kickoff date (column D)
# untouched, exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-FixtureRows($rowA, $rowB) {
    # Column B (match id) - single cell swap
    $idA = $ws.Range("B$rowA").Value2
    $idB = $ws.Range("B$rowB").Value2
    $ws.Range("B$rowA").Value2 = $idB
    $ws.Range("B$rowB").Value2 = $idA

    # Columns E through AD (teams, score, result, odds) - block swap
    $blockA = $ws.Range("E$rowA`:AD$rowA").Value2
    $blockB = $ws.Range("E$rowB`:AD$rowB").Value2
    $ws.Range("E$rowA`:AD$rowA").Value2 = $blockB
    $ws.Range("E$rowB`:AD$rowB").Value2 = $blockA
}

# Rows 4 and 5: JK Tallinna Kalev vs JK Trans Narva / JK Tammeka Tartu vs Harju JK Laagri
Swap-FixtureRows 4 5

# Rows 115 and 116: FC Kuressaare vs FC Levadia Tallinn / JK Nomme Kalju vs JK Trans Narva
Swap-FixtureRows 115 116
